# Applies the "Alvearie -> LinuxForHealth" rebrand edit described by the
# commit "Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@..."
#
# Sheet "Metadata": update URL, Version, Date and Publisher values.
# Sheet "Elements": clear the stray Constraint(s) text that had leaked into
# the "Extension" summary row, and let the Fixed Value cell cascade from the
# URL update above.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------
$wsMetadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-cost-scale"
$wsMetadata.Range("B3").Value = "8.0.0"
$wsMetadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMetadata.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------------
# Row 2 ("Extension") had an erroneous Constraint(s) value that belongs to
# the "Extension.extension" row only; clear it here.
$wsElements.Range("AI2").Value = ""

# Row 5's Fixed Value mirrors the StructureDefinition URL (shared string),
# so it must be refreshed explicitly to stay in sync with the new URL.
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-cost-scale"
